$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 620.27026
$ws.Range("I33").Value = 143.26666
$ws.Range("J33").Value = 2664.5715
$ws.Range("K33").Value = 143.26666
$ws.Range("L33").Value = 2664.5715
$ws.Range("M33").Value = 85.73334
$ws.Range("N33").Value = -3122.5715

# Row 96
$ws.Range("H96").Value = 1999
$ws.Range("I96").Value = 998.6667
$ws.Range("J96").Value = 5000
$ws.Range("K96").Value = 2996.0001
$ws.Range("L96").Value = 15000
$ws.Range("M96").Value = -1623.0001
$ws.Range("N96").Value = -17746

# Row 125
$ws.Range("H125").Value = 3609.7
$ws.Range("I125").Value = 5516
$ws.Range("J125").Value = 3133.125
$ws.Range("K125").Value = 49644
$ws.Range("L125").Value = 28198.125
$ws.Range("M125").Value = -47184
$ws.Range("N125").Value = -33118.125

# Row 132
$ws.Range("H132").Value = 5560917
$ws.Range("I132").Value = 6103138
$ws.Range("J132").Value = 3150
$ws.Range("K132").Value = 18309414
$ws.Range("L132").Value = 9450
$ws.Range("M132").Value = -18306884
$ws.Range("N132").Value = -14510

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 25271.213
$ws.Range("I32").Value = 5272.2114
$ws.Range("K32").Value = 5272.2114
$ws.Range("M32").Value = -4985.2114

# Row 74
$ws.Range("H74").Value = 988.3182
$ws.Range("I74").Value = 961.58826
$ws.Range("J74").Value = 1079.2
$ws.Range("K74").Value = 961.58826
$ws.Range("L74").Value = 1079.2
$ws.Range("M74").Value = -87.58825999999999
$ws.Range("N74").Value = -2827.2

# Row 77
$ws.Range("H77").Value = 988.3182
$ws.Range("I77").Value = 961.58826
$ws.Range("J77").Value = 1079.2
$ws.Range("K77").Value = 4807.9413
$ws.Range("L77").Value = 5396
$ws.Range("M77").Value = -439.9413000000004
$ws.Range("N77").Value = -14132

# Row 96
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

# Row 102
$ws.Range("H102").Value = 51756.15
$ws.Range("J102").Value = 1811.0834
$ws.Range("L102").Value = 1811.0834
$ws.Range("N102").Value = -5055.0834

$ws = $wb.Worksheets.Item("CRP")
# Row 36
$ws.Range("H36").Value = 5000
$ws.Range("I36").Value = 5000
$ws.Range("K36").Value = 5000
$ws.Range("M36").Value = -4612

# Row 40
$ws.Range("H40").Value = 5000
$ws.Range("I40").Value = 5000
$ws.Range("K40").Value = 5000
$ws.Range("M40").Value = -4840

# Row 74
$ws.Range("H74").Value = 25823.182
$ws.Range("J74").Value = 25823.182
$ws.Range("L74").Value = 25823.182
$ws.Range("N74").Value = -27571.182

# Row 77
$ws.Range("H77").Value = 25823.182
$ws.Range("J77").Value = 25823.182
$ws.Range("L77").Value = 77469.546
$ws.Range("N77").Value = -86205.546

# Row 88
$ws.Range("H88").Value = 23134.4
$ws.Range("J88").Value = 23134.4
$ws.Range("L88").Value = 23134.4
$ws.Range("N88").Value = -23946.4

# Row 91
$ws.Range("H91").Value = 23134.4
$ws.Range("J91").Value = 23134.4
$ws.Range("L91").Value = 23134.4
$ws.Range("N91").Value = -25942.4

# Row 132
$ws.Range("H132").Value = 4464.654
$ws.Range("I132").Value = 4189.1577
$ws.Range("J132").Value = 5212.4287
$ws.Range("K132").Value = 12567.4731
$ws.Range("L132").Value = 15637.2861
$ws.Range("M132").Value = -10037.4731
$ws.Range("N132").Value = -20697.2861

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 1722.3077
$ws.Range("I4").Value = 195
$ws.Range("J4").Value = 2000
$ws.Range("K4").Value = 585
$ws.Range("L4").Value = 6000
$ws.Range("M4").Value = -473
$ws.Range("N4").Value = -6224

# Row 9
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

# Row 86
$ws.Range("H86").Value = 20000586
$ws.Range("J86").Value = 25000616
$ws.Range("L86").Value = 75001848
$ws.Range("N86").Value = -75004220

# Row 89
$ws.Range("H89").Value = 20000586
$ws.Range("J89").Value = 25000616
$ws.Range("L89").Value = 225005544
$ws.Range("N89").Value = -225017400

# Row 113
$ws.Range("H113").Value = 966.6818
$ws.Range("I113").Value = 1699.5
$ws.Range("K113").Value = 5098.5
$ws.Range("M113").Value = -2928.5

# Row 122
$ws.Range("H122").Value = 554
$ws.Range("I122").Value = 554
$ws.Range("K122").Value = 4986
$ws.Range("M122").Value = -2536

# Row 132
$ws.Range("H132").Value = 1790.9259
$ws.Range("I132").Value = 772.2222
$ws.Range("J132").Value = 2300.2778
$ws.Range("K132").Value = 6949.999800000001
$ws.Range("L132").Value = 20702.5002
$ws.Range("M132").Value = -4419.999800000001
$ws.Range("N132").Value = -25762.5002

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 1033.091
$ws.Range("I122").Value = 819
$ws.Range("J122").Value = 1113.375
$ws.Range("K122").Value = 2457
$ws.Range("L122").Value = 3340.125
$ws.Range("M122").Value = -7
$ws.Range("N122").Value = -8240.125

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1477.4445
$ws.Range("I22").Value = 1490.8572
$ws.Range("J22").Value = 1463
$ws.Range("K22").Value = 1490.8572
$ws.Range("L22").Value = 1463
$ws.Range("M22").Value = -1195.8572
$ws.Range("N22").Value = -2053

# Row 27
$ws.Range("H27").Value = 1477.4445
$ws.Range("I27").Value = 1490.8572
$ws.Range("J27").Value = 1463
$ws.Range("K27").Value = 1490.8572
$ws.Range("L27").Value = 1463
$ws.Range("M27").Value = -1383.8572
$ws.Range("N27").Value = -1677

# Row 46
$ws.Range("H46").Value = 506450
$ws.Range("I46").Value = 344.36365
$ws.Range("J46").Value = 1125023.5
$ws.Range("K46").Value = 344.36365
$ws.Range("L46").Value = 1125023.5
$ws.Range("M46").Value = -156.36365
$ws.Range("N46").Value = -1125399.5

# Row 136
$ws.Range("H136").Value = 2364.7058
$ws.Range("I136").Value = 2228.5715
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 6685.7145
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -4135.7145
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 1346.35
$ws.Range("I126").Value = 995.4375
$ws.Range("J126").Value = 2750
$ws.Range("K126").Value = 2986.3125
$ws.Range("L126").Value = 8250
$ws.Range("M126").Value = -516.3125
$ws.Range("N126").Value = -13190

# Row 136
$ws.Range("H136").Value = 1484.6061
$ws.Range("I136").Value = 580.2258
$ws.Range("J136").Value = 15502.5
$ws.Range("K136").Value = 1740.6774
$ws.Range("L136").Value = 46507.5
$ws.Range("M136").Value = 809.3226
$ws.Range("N136").Value = -51607.5

Write-Output "Applied all cell updates"